$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 49
$ws.Range("N49").ClearContents()
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("H49").Value = 199
# Row 62
$ws.Range("M62").Value = -5274
$ws.Range("H62").Value = 5952
$ws.Range("N62").Value = -7254
$ws.Range("I62").Value = 5898
$ws.Range("J62").Value = 6006
$ws.Range("K62").Value = 5898
$ws.Range("L62").Value = 6006
# Row 65
$ws.Range("M65").Value = -26370
$ws.Range("H65").Value = 5952
$ws.Range("L65").Value = 30030
$ws.Range("J65").Value = 6006
$ws.Range("I65").Value = 5898
$ws.Range("N65").Value = -36270
$ws.Range("K65").Value = 29490
# Row 86
$ws.Range("J86").Value = 4000
$ws.Range("M86").Value = -1745.3333
$ws.Range("I86").Value = 2868.3333
$ws.Range("K86").Value = 2868.3333
$ws.Range("H86").Value = 3151.25
$ws.Range("N86").Value = -6246
$ws.Range("L86").Value = 4000
# Row 89
$ws.Range("K89").Value = 14341.6665
$ws.Range("J89").Value = 4000
$ws.Range("L89").Value = 20000
$ws.Range("I89").Value = 2868.3333
$ws.Range("N89").Value = -31232
$ws.Range("M89").Value = -8725.666499999999
$ws.Range("H89").Value = 3151.25
# Row 132
$ws.Range("I132").Value = 2322.8
$ws.Range("K132").Value = 6968.400000000001
$ws.Range("M132").Value = -4438.400000000001
$ws.Range("H132").Value = 2322.8
# Row 137
$ws.Range("H137").Value = 2628
$ws.Range("I137").Value = 1293.5
$ws.Range("M137").Value = -1330.5
$ws.Range("K137").Value = 3880.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3550
$ws.Range("K2").Value = 3550
$ws.Range("M2").Value = -3437
$ws.Range("I2").Value = 3550
# Row 32
$ws.Range("H32").Value = 2795.2334
$ws.Range("M32").Value = -2438.3845
$ws.Range("K32").Value = 2725.3845
$ws.Range("I32").Value = 2725.3845
# Row 88
$ws.Range("L88").Value = 4999
$ws.Range("I88").Value = 5000
$ws.Range("M88").Value = -4594
$ws.Range("J88").Value = 4999
$ws.Range("K88").Value = 5000
$ws.Range("N88").Value = -5811
$ws.Range("H88").Value = 4999.5
# Row 91
$ws.Range("L91").Value = 4999
$ws.Range("N91").Value = -7807
$ws.Range("J91").Value = 4999
$ws.Range("H91").Value = 4999.5
$ws.Range("M91").Value = -3596
$ws.Range("I91").Value = 5000
$ws.Range("K91").Value = 5000
# Row 116
$ws.Range("I116").Value = 3550
$ws.Range("H116").Value = 3550
$ws.Range("K116").Value = 3550
$ws.Range("M116").Value = -1256
# Row 122
$ws.Range("L122").Value = 3163.5
$ws.Range("H122").Value = 1569.8334
$ws.Range("N122").Value = -8063.5
$ws.Range("J122").Value = 1054.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3550
$ws.Range("K3").Value = 3550
$ws.Range("I3").Value = 3550
$ws.Range("M3").Value = -3436
# Row 99
$ws.Range("I99").Value = 997.25
$ws.Range("H99").Value = 997.25
$ws.Range("M99").Value = 500.75
$ws.Range("K99").Value = 997.25
# Row 107
$ws.Range("I107").Value = 2332.6667
$ws.Range("M107").Value = -412.6667000000002
$ws.Range("K107").Value = 2332.6667
$ws.Range("H107").Value = 2856.5715

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("I99").Value = 3599
$ws.Range("H99").Value = 3599.5
$ws.Range("J99").Value = 3600
$ws.Range("M99").Value = -2101
$ws.Range("N99").Value = -6596
$ws.Range("K99").Value = 3599
$ws.Range("L99").Value = 3600
# Row 126
$ws.Range("L126").Value = 10800
$ws.Range("H126").Value = 3599.5
$ws.Range("M126").Value = -8327
$ws.Range("J126").Value = 3600
$ws.Range("K126").Value = 10797
$ws.Range("I126").Value = 3599
$ws.Range("N126").Value = -15740

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 102
$ws.Range("L102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("H102").Value = 550
$ws.Range("N102").ClearContents()
# Row 113
$ws.Range("H113").Value = 958.6
$ws.Range("K113").Value = 2846.25
$ws.Range("M113").Value = -676.25
$ws.Range("I113").Value = 948.75

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("J80").Value = 7251
$ws.Range("H80").Value = 6287
$ws.Range("N80").Value = -9247
$ws.Range("L80").Value = 7251
$ws.Range("I80").Value = 5001.6665
$ws.Range("K80").Value = 5001.6665
$ws.Range("M80").Value = -4003.6665
# Row 83
$ws.Range("M83").Value = -20016.3325
$ws.Range("L83").Value = 36255
$ws.Range("K83").Value = 25008.3325
$ws.Range("I83").Value = 5001.6665
$ws.Range("N83").Value = -46239
$ws.Range("H83").Value = 6287
$ws.Range("J83").Value = 7251
# Row 107
$ws.Range("I107").Value = 450
$ws.Range("M107").Value = 1470
$ws.Range("K107").Value = 450
$ws.Range("H107").Value = 450

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2072.8572
$ws.Range("J61").Value = 2071.5
$ws.Range("N61").Value = -2475.5
$ws.Range("K61").Value = 2074.6667
$ws.Range("L61").Value = 2071.5
$ws.Range("M61").Value = -1872.6667
$ws.Range("I61").Value = 2074.6667
# Row 82
$ws.Range("N82").Value = -1941.25
$ws.Range("K82").Value = 1899.75
$ws.Range("I82").Value = 1899.75
$ws.Range("H82").Value = 1559.5
$ws.Range("J82").Value = 1219.25
$ws.Range("L82").Value = 1219.25
$ws.Range("M82").Value = -1538.75
# Row 85
$ws.Range("H85").Value = 1559.5
$ws.Range("K85").Value = 1899.75
$ws.Range("N85").Value = -3715.25
$ws.Range("I85").Value = 1899.75
$ws.Range("J85").Value = 1219.25
$ws.Range("M85").Value = -651.75
$ws.Range("L85").Value = 1219.25
# Row 100
$ws.Range("M100").Value = -2958.5
$ws.Range("H100").Value = 3499.5
$ws.Range("I100").Value = 3499.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3499.5
$ws.Range("N100").ClearContents()
$ws.Range("L100").Value = 0
# Row 113
$ws.Range("H113").Value = 2072.8572
$ws.Range("K113").Value = 2074.6667
$ws.Range("N113").Value = -6411.5
$ws.Range("J113").Value = 2071.5
$ws.Range("M113").Value = 95.33329999999978
$ws.Range("I113").Value = 2074.6667
$ws.Range("L113").Value = 2071.5

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 499.33334
$ws.Range("K113").Value = 1498.00002
$ws.Range("M113").Value = 671.9999800000001
$ws.Range("I113").Value = 499.33334
